$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Step 1: duplicate the "2021-Q4" sheet to become the new "2022-Q1" sheet ---
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($null, $templateSheet)
$ws = $wb.Worksheets.Item(3)
$ws.Name = "2022-Q1"

# --- Step 2: extend formatting for the additional rows (10-21) by copying row 9 formats down ---
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H21").PasteSpecial(-4122)

# --- Step 3: write header row ---
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# --- Step 4: write data rows 2-21 ---
$row = 2
$ws.Cells.Item($row,1).Value = 0
Set-TextCell $ws.Cells.Item($row,2) "010336"
Set-TextCell $ws.Cells.Item($row,3) "中欧悦享生活混合A"
Set-TextCell $ws.Cells.Item($row,4) "44.28"
Set-TextCell $ws.Cells.Item($row,5) "90.44"
Set-TextCell $ws.Cells.Item($row,6) "4.64"
Set-TextCell $ws.Cells.Item($row,7) "2.0546"
$ws.Cells.Item($row,8).Value = 8

$row = 3
$ws.Cells.Item($row,1).Value = 1
Set-TextCell $ws.Cells.Item($row,2) "006529"
Set-TextCell $ws.Cells.Item($row,3) "中欧匠心两年持有期混合A"
Set-TextCell $ws.Cells.Item($row,4) "43.69"
Set-TextCell $ws.Cells.Item($row,5) "88.30"
Set-TextCell $ws.Cells.Item($row,6) "3.03"
Set-TextCell $ws.Cells.Item($row,7) "1.3238"
$ws.Cells.Item($row,8).Value = 8

$row = 4
$ws.Cells.Item($row,1).Value = 2
Set-TextCell $ws.Cells.Item($row,2) "007291"
Set-TextCell $ws.Cells.Item($row,3) "汇丰晋信港股通双核策略混合"
Set-TextCell $ws.Cells.Item($row,4) "7.80"
Set-TextCell $ws.Cells.Item($row,5) "92.66"
Set-TextCell $ws.Cells.Item($row,6) "4.48"
Set-TextCell $ws.Cells.Item($row,7) "0.3494"
$ws.Cells.Item($row,8).Value = 8

$row = 5
$ws.Cells.Item($row,1).Value = 3
Set-TextCell $ws.Cells.Item($row,2) "002332"
Set-TextCell $ws.Cells.Item($row,3) "汇丰晋信沪港深股票A"
Set-TextCell $ws.Cells.Item($row,4) "7.90"
Set-TextCell $ws.Cells.Item($row,5) "92.60"
Set-TextCell $ws.Cells.Item($row,6) "3.87"
Set-TextCell $ws.Cells.Item($row,7) "0.3057"
$ws.Cells.Item($row,8).Value = 10

$row = 6
$ws.Cells.Item($row,1).Value = 4
Set-TextCell $ws.Cells.Item($row,2) "012744"
Set-TextCell $ws.Cells.Item($row,3) "光大保德信品质生活混合型证券投资基金A"
Set-TextCell $ws.Cells.Item($row,4) "6.91"
Set-TextCell $ws.Cells.Item($row,5) "84.96"
Set-TextCell $ws.Cells.Item($row,6) "3.49"
Set-TextCell $ws.Cells.Item($row,7) "0.2412"
$ws.Cells.Item($row,8).Value = 9

$row = 7
$ws.Cells.Item($row,1).Value = 5
Set-TextCell $ws.Cells.Item($row,2) "005620"
Set-TextCell $ws.Cells.Item($row,3) "中欧品质消费股票A"
Set-TextCell $ws.Cells.Item($row,4) "3.74"
Set-TextCell $ws.Cells.Item($row,5) "90.47"
Set-TextCell $ws.Cells.Item($row,6) "5.65"
Set-TextCell $ws.Cells.Item($row,7) "0.2113"
$ws.Cells.Item($row,8).Value = 8

$row = 8
$ws.Cells.Item($row,1).Value = 6
Set-TextCell $ws.Cells.Item($row,2) "870017"
Set-TextCell $ws.Cells.Item($row,3) "广发资管消费精选灵活配置混合"
Set-TextCell $ws.Cells.Item($row,4) "2.56"
Set-TextCell $ws.Cells.Item($row,5) "93.50"
Set-TextCell $ws.Cells.Item($row,6) "5.93"
Set-TextCell $ws.Cells.Item($row,7) "0.1518"
$ws.Cells.Item($row,8).Value = 9

$row = 9
$ws.Cells.Item($row,1).Value = 7
Set-TextCell $ws.Cells.Item($row,2) "006530"
Set-TextCell $ws.Cells.Item($row,3) "中欧匠心两年持有期混合C"
Set-TextCell $ws.Cells.Item($row,4) "4.59"
Set-TextCell $ws.Cells.Item($row,5) "88.30"
Set-TextCell $ws.Cells.Item($row,6) "3.03"
Set-TextCell $ws.Cells.Item($row,7) "0.1391"
$ws.Cells.Item($row,8).Value = 8

$row = 10
$ws.Cells.Item($row,1).Value = 8
Set-TextCell $ws.Cells.Item($row,2) "009877"
Set-TextCell $ws.Cells.Item($row,3) "中银内核驱动股票"
Set-TextCell $ws.Cells.Item($row,4) "2.82"
Set-TextCell $ws.Cells.Item($row,5) "83.99"
Set-TextCell $ws.Cells.Item($row,6) "3.96"
Set-TextCell $ws.Cells.Item($row,7) "0.1117"
$ws.Cells.Item($row,8).Value = 10

$row = 11
$ws.Cells.Item($row,1).Value = 9
Set-TextCell $ws.Cells.Item($row,2) "005621"
Set-TextCell $ws.Cells.Item($row,3) "中欧品质消费股票C"
Set-TextCell $ws.Cells.Item($row,4) "1.11"
Set-TextCell $ws.Cells.Item($row,5) "90.47"
Set-TextCell $ws.Cells.Item($row,6) "5.65"
Set-TextCell $ws.Cells.Item($row,7) "0.0627"
$ws.Cells.Item($row,8).Value = 8

$row = 12
$ws.Cells.Item($row,1).Value = 10
Set-TextCell $ws.Cells.Item($row,2) "007109"
Set-TextCell $ws.Cells.Item($row,3) "南方沪港深核心优势混合"
Set-TextCell $ws.Cells.Item($row,4) "1.82"
Set-TextCell $ws.Cells.Item($row,5) "87.54"
Set-TextCell $ws.Cells.Item($row,6) "3.03"
Set-TextCell $ws.Cells.Item($row,7) "0.0551"
$ws.Cells.Item($row,8).Value = 9

$row = 13
$ws.Cells.Item($row,1).Value = 11
Set-TextCell $ws.Cells.Item($row,2) "010337"
Set-TextCell $ws.Cells.Item($row,3) "中欧悦享生活混合C"
Set-TextCell $ws.Cells.Item($row,4) "1.08"
Set-TextCell $ws.Cells.Item($row,5) "90.44"
Set-TextCell $ws.Cells.Item($row,6) "4.64"
Set-TextCell $ws.Cells.Item($row,7) "0.0501"
$ws.Cells.Item($row,8).Value = 8

$row = 14
$ws.Cells.Item($row,1).Value = 12
Set-TextCell $ws.Cells.Item($row,2) "002333"
Set-TextCell $ws.Cells.Item($row,3) "汇丰晋信沪港深股票C"
Set-TextCell $ws.Cells.Item($row,4) "1.23"
Set-TextCell $ws.Cells.Item($row,5) "92.60"
Set-TextCell $ws.Cells.Item($row,6) "3.87"
Set-TextCell $ws.Cells.Item($row,7) "0.0476"
$ws.Cells.Item($row,8).Value = 10

$row = 15
$ws.Cells.Item($row,1).Value = 13
Set-TextCell $ws.Cells.Item($row,2) "241001"
Set-TextCell $ws.Cells.Item($row,3) "华宝海外中国混合(QDII)"
Set-TextCell $ws.Cells.Item($row,4) "0.83"
Set-TextCell $ws.Cells.Item($row,5) "86.89"
Set-TextCell $ws.Cells.Item($row,6) "4.13"
Set-TextCell $ws.Cells.Item($row,7) "0.0343"
$ws.Cells.Item($row,8).Value = 9

$row = 16
$ws.Cells.Item($row,1).Value = 14
Set-TextCell $ws.Cells.Item($row,2) "004099"
Set-TextCell $ws.Cells.Item($row,3) "前海开源沪港深景气行业精选灵活配置混合"
Set-TextCell $ws.Cells.Item($row,4) "0.41"
Set-TextCell $ws.Cells.Item($row,5) "93.07"
Set-TextCell $ws.Cells.Item($row,6) "8.20"
Set-TextCell $ws.Cells.Item($row,7) "0.0336"
$ws.Cells.Item($row,8).Value = 8

$row = 17
$ws.Cells.Item($row,1).Value = 15
Set-TextCell $ws.Cells.Item($row,2) "012315"
Set-TextCell $ws.Cells.Item($row,3) "创金合信港股通成长股票型发起式证券投资基金A"
Set-TextCell $ws.Cells.Item($row,4) "0.19"
Set-TextCell $ws.Cells.Item($row,5) "83.49"
Set-TextCell $ws.Cells.Item($row,6) "5.68"
Set-TextCell $ws.Cells.Item($row,7) "0.0108"
$ws.Cells.Item($row,8).Value = 8

$row = 18
$ws.Cells.Item($row,1).Value = 16
Set-TextCell $ws.Cells.Item($row,2) "012758"
Set-TextCell $ws.Cells.Item($row,3) "光大保德信品质生活混合型证券投资基金C"
Set-TextCell $ws.Cells.Item($row,4) "0.31"
Set-TextCell $ws.Cells.Item($row,5) "84.96"
Set-TextCell $ws.Cells.Item($row,6) "3.49"
Set-TextCell $ws.Cells.Item($row,7) "0.0108"
$ws.Cells.Item($row,8).Value = 9

$row = 19
$ws.Cells.Item($row,1).Value = 17
Set-TextCell $ws.Cells.Item($row,2) "009733"
Set-TextCell $ws.Cells.Item($row,3) "创金合信港股通大消费精选股票A"
Set-TextCell $ws.Cells.Item($row,4) "0.13"
Set-TextCell $ws.Cells.Item($row,5) "82.28"
Set-TextCell $ws.Cells.Item($row,6) "4.51"
Set-TextCell $ws.Cells.Item($row,7) "0.0059"
$ws.Cells.Item($row,8).Value = 7

$row = 20
$ws.Cells.Item($row,1).Value = 18
Set-TextCell $ws.Cells.Item($row,2) "012316"
Set-TextCell $ws.Cells.Item($row,3) "创金合信港股通成长股票型发起式证券投资基金C"
Set-TextCell $ws.Cells.Item($row,4) "0.10"
Set-TextCell $ws.Cells.Item($row,5) "83.49"
Set-TextCell $ws.Cells.Item($row,6) "5.68"
Set-TextCell $ws.Cells.Item($row,7) "0.0057"
$ws.Cells.Item($row,8).Value = 8

$row = 21
$ws.Cells.Item($row,1).Value = 19
Set-TextCell $ws.Cells.Item($row,2) "009734"
Set-TextCell $ws.Cells.Item($row,3) "创金合信港股通大消费精选股票C"
Set-TextCell $ws.Cells.Item($row,4) "0.07"
Set-TextCell $ws.Cells.Item($row,5) "82.28"
Set-TextCell $ws.Cells.Item($row,6) "4.51"
Set-TextCell $ws.Cells.Item($row,7) "0.0032"
$ws.Cells.Item($row,8).Value = 7

# --- Step 5: update the "总计" (summary) sheet: insert a new row for 2022-Q1 at the top of the data ---
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()

# Re-apply formatting that Insert() leaves inconsistent: column A should carry the bordered/bold style,
# columns B:D should have the plain (unstyled) format used by the other data rows.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)
$zj.Range("B3:D3").Copy()
$zj.Range("B2:D2").PasteSpecial(-4122)

$zj.Cells.Item(2,1).Value = 0
Set-TextCell $zj.Cells.Item(2,2) "2022-Q1"
$zj.Cells.Item(2,3).Value = 20
$zj.Cells.Item(2,4).Value = 5.21

# Renumber the index column for the rows that got pushed down
$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(4,1).Value = 2

Write-Host "Edit complete"
